# Agregada captura de periodicidad IVA
#
# Inserts a new "Periodicidad de la declaración" / "Cuatrimestral" row right
# after the existing "Período" row (old row 2), pushing every subsequent
# Concepto/Campo/Valor/Año/Razón Social row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the old row 3 ("Ingresos por operaciones gravadas
# al 5%"), shifting it (and everything below it) down to row 4.
$ws.Rows.Item(3).Insert()

# Match the row height/formatting convention used by every other data row.
$ws.Rows.Item(3).RowHeight = 36

# Populate the newly inserted row 3 with the periodicidad fields.
$ws.Range("A3").Value = "Periodicidad de la declaración"
$ws.Range("B3").Value = 24
$ws.Range("C3").Value = "Cuatrimestral"
$ws.Range("D3").Value = 2019
$ws.Range("E3").Value = "INVERSIONES ORTIZ VASQUEZ HERMANOS S A S"
